$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.149.13'
$ws.Range('E2').Value = '  -3.23%  '
$ws.Range('D3').Value = '1.849.45'
$ws.Range('E3').Value = '  -2.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7029'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '238.49'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.77%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3055'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07471'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.39'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08123'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Value = '1.869.37'
$ws.Range('E12').Value = '  -2.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7256'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.224'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '88.73'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.78%  '
$ws.Range('D16').Value = '29.243.24'
$ws.Range('E16').Value = '  -2.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.757'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.73%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '238.18'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.08'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007631'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.0000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').Value = '2.114.47'
$ws.Range('E22').Value = '  -2.11%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.576'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.000'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1452'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.07'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.977'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.398'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.539'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.494'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.975'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05171'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.187'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.66%  '
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.044'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.71%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6997'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -9.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.652'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01866'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.94%  '
$ws.Range('E40').Value = '  -3.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9408'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.001'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.55%  '
$ws.Range('D43').Value = '1.080.40'
$ws.Range('E43').Value = '  -1.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4285'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '70.14'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9998'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.50'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.89%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '2.008.13'
$ws.Range('E48').Value = '  -3.08%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.744'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.035'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.144'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.81%  '
